# Edit the "Controls" sheet: update a couple of parameter values and
# move the active selection, as part of documenting the simulate_data
# function's input parameters.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controls")
$ws.Activate()

# n_sims: 100 -> 500
$ws.Range("B2").Value = 500

# n_srv_fleets: 2 -> 1
$ws.Range("B6").Value = 1

# Move the active selection/cell to B3
$ws.Range("B3").Select()
